$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title: "Docker 101: Run your first Docker container today" (+ trailing
#    space run) -> "Docker 101: Run Your First Docker Container Today "
#    with the font shrunk from 44 -> 42 half-points, and the document's
#    "_GoBack" bookmark relocated into the middle of the new title text
#    (GHC's capitalization scheme / last-edit location).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Find.Execute("Run your first Docker container today", $true, $false, $false, $false, $false, $true, 1, $false, "Run Your First Docker Container Today", 2)

# Shrink the title font (both the Western and complex-script sizes).
$titlePara.Range.Font.Size = 21
$titlePara.Range.Font.SizeBi = 21

# Create a temporary hard boundary right before "oday" so that the upcoming
# bookmark insertion/relocation doesn't cause that tail to re-merge with the
# "...Container T" run.
$splitProbe = $titlePara.Range.Duplicate
$splitProbe.Find.Execute("oday", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailSplitPoint = $d.Range($splitProbe.Start, $splitProbe.Start)
$d.Bookmarks.Add("zzTmpSplit", $tailSplitPoint)

# Remove the bookmark from its old home (inside the "will be available"
# paragraph further down) before re-adding it at the new location - Word
# only ever has one "_GoBack".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Place "_GoBack" right after "...Run Your Fi" (splitting the title run).
$headProbe = $titlePara.Range.Duplicate
$headProbe.Find.Execute("rst Docker", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headSplitPoint = $d.Range($headProbe.Start, $headProbe.Start)
$d.Bookmarks.Add("_GoBack", $headSplitPoint)

# Drop the temporary boundary now that the real run split has materialized.
$d.Bookmarks("zzTmpSplit").Delete()

# ---------------------------------------------------------------------------
# 2) "... accounts will be available after the workshop ..." previously had
#    "_GoBack" splitting " will " from "be available after the workshop";
#    now that the bookmark moved into the title, those two runs collapse
#    back into one.
# ---------------------------------------------------------------------------
$body = $d.Content

$mergeProbe = $body.Duplicate
$mergeProbe.Find.Execute(" will be available after the workshop", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStart = $d.Range($mergeProbe.Start, $mergeProbe.Start)
$d.Bookmarks.Add("zzMergeStart", $mergeStart)

$stopProbe = $body.Duplicate
$stopProbe.Find.Execute("be available after the workshop", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStop = $d.Range($stopProbe.End, $stopProbe.End)
$d.Bookmarks.Add("zzMergeStop", $mergeStop)

$mergeRange = $body.Duplicate
$mergeRange.Find.Execute(" will be available after the workshop", $true, $false, $false, $false, $false, $true, 1, $false, " will be available after the workshop", 2)

$d.Bookmarks("zzMergeStart").Delete()
$d.Bookmarks("zzMergeStop").Delete()

Write-Host "Title now: [$($titlePara.Range.Text)]"
